# "updated program so that overlay of spectra works correctly"
#
# This adds a new "NOESY" spectrum sheet (mirroring the existing HMBC-style
# 2D-peak-list sheets: f2/f1/Width f2/Width f1/Volume/Type/Flags/
# Impurity-Compound/Annotation header row) after the HMBC sheet, makes it
# the active sheet/tab, refreshes the COSY sheet's selected overlay range,
# and tidies up the molecule sheet's column widths so the SMILES column is
# fully visible.

$wb = $excel.ActiveWorkbook

# --- COSY: re-select the full header row (the overlay range used when the
#     spectra are plotted together) -----------------------------------
$cosy = $wb.Worksheets.Item("COSY")
$cosy.Range("A1:K1").Select()

# --- molecule: widen the Name/SMILES columns so the text isn't truncated
$mol = $wb.Worksheets.Item("molecule")
$mol.Columns.Item(2).ColumnWidth = 9.5
$mol.Columns.Item(3).ColumnWidth = 20.666666666666668

# --- add the new NOESY sheet after the last existing sheet (HMBC) -------
$hmbc = $wb.Worksheets.Item("HMBC")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$noesy = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$noesy.Name = "NOESY"

# Bring over the standard 2D peak-list header (f2 (ppm), f1 (ppm),
# Width f2, Width f1, Volume, Type, Flags, Impurity/Compound, Annotation)
# along with its formatting, same as COSY/HSQC/HMBC use.
$hmbc.Range("B1:K1").Copy($noesy.Range("B1"))

# --- H1_1D: rows no longer need the extra-tall 21pt override ------------
$h1 = $wb.Worksheets.Item("H1_1D")
$h1.Range("A2:I6").EntireRow.AutoFit()

# --- finally, select the header row on NOESY and make it the active tab
$noesy.Range("A1:K1").Select()
$noesy.Activate()
